$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update BLEU score (B11)
$ws.Range("B11").Value = 0.1879823945296794

# Update Code BLEU (B12)
$ws.Range("B12").Value = 0.3381137531680278

# Update the Code BLEU detail note (C12)
$ws.Range("C12").Value = "{'codebleu': 0.3381137531680278, 'ngram_match_score': 0.18798239452967938, 'weighted_ngram_match_score': 0.1974289091632388, 'syntax_match_score': 0.5584415584415584, 'dataflow_match_score': 0.40860215053763443}"

# Update Embeddings and Cosine similarity (B13)
$ws.Range("B13").Value = 0.9048691301042167
